$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("runs")

$ws.Range("B1").Value = 3
$ws.Range("B2").Value = 3

$ws.Activate()
$ws.Range("B2").Select()
